$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three worker rows that are no longer part of this statement
# (EDER DE JESUS AVILA BERRIO, JUAN FELIPE GUERRERO LONDOÑO, JORGE RICARDO GUERRERO MOLANO).
# Delete bottom-up so earlier row numbers stay valid while we work.
$ws.Rows.Item(19).EntireRow.Delete()
$ws.Rows.Item(18).EntireRow.Delete()
$ws.Rows.Item(16).EntireRow.Delete()

# Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 56934
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3

# The two remaining "MARIA ELENA PEÑA CHAMORRO" rows (now rows 17 and 18) swap which
# overdue period they report, and both now carry the updated base salary.
$ws.Range("E17").Value = "2112"
$ws.Range("F17").Value = 3634
$ws.Range("G17").Value = 908526

$ws.Range("E18").Value = "2201"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526
